$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 29800
$ws.Range("J44").Value = 29800
$ws.Range("L44").Value = 29800
$ws.Range("N44").Value = -30724
$ws.Range("H47").Value = 24900
$ws.Range("I47").Value = 20000
$ws.Range("J47").Value = 29800
$ws.Range("K47").Value = 20000
$ws.Range("L47").Value = 29800
$ws.Range("M47").Value = -19028
$ws.Range("N47").Value = -31744
$ws.Range("H74").Value = 6717.207
$ws.Range("I74").Value = 12650
$ws.Range("J74").Value = 3594.6843
$ws.Range("K74").Value = 12650
$ws.Range("L74").Value = 3594.6843
$ws.Range("M74").Value = -11714
$ws.Range("N74").Value = -5466.6843
$ws.Range("H77").Value = 6717.207
$ws.Range("I77").Value = 12650
$ws.Range("J77").Value = 3594.6843
$ws.Range("K77").Value = 63250
$ws.Range("L77").Value = 17973.4215
$ws.Range("M77").Value = -58570
$ws.Range("N77").Value = -27333.4215
$ws.Range("H80").Value = 1213
$ws.Range("J80").Value = 1775
$ws.Range("L80").Value = 5325
$ws.Range("N80").Value = -7321
$ws.Range("H83").Value = 1213
$ws.Range("J83").Value = 1775
$ws.Range("L83").Value = 15975
$ws.Range("N83").Value = -25959
$ws.Range("H100").Value = 45456030
$ws.Range("I100").Value = 1634.375
$ws.Range("K100").Value = 1634.375
$ws.Range("M100").Value = -1093.375
$ws.Range("H103").Value = 6061209
$ws.Range("I103").Value = 481.66666
$ws.Range("J103").Value = 16667482
$ws.Range("K103").Value = 1444.99998
$ws.Range("L103").Value = 50002446
$ws.Range("M103").Value = -858.9999800000001
$ws.Range("N103").Value = -50003618

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2220.51
$ws.Range("I32").Value = 2220.51
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2220.51
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1933.51
$ws.Range("N32").ClearContents()
$ws.Range("H102").Value = 71430340
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 83335060
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 83335060
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -83338304
$ws.Range("H122").Value = 21064.727
$ws.Range("I122").Value = 30876
$ws.Range("K122").Value = 92628
$ws.Range("M122").Value = -90178

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1681.6666
$ws.Range("I99").Value = 1097.5
$ws.Range("J99").Value = 2850
$ws.Range("K99").Value = 1097.5
$ws.Range("L99").Value = 2850
$ws.Range("M99").Value = 400.5
$ws.Range("N99").Value = -5846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3535.6902
$ws.Range("I31").Value = 2035.7142
$ws.Range("J31").Value = 4165.68
$ws.Range("K31").Value = 2035.7142
$ws.Range("L31").Value = 4165.68
$ws.Range("M31").Value = -1740.7142
$ws.Range("N31").Value = -4755.68
$ws.Range("H34").Value = 3535.6902
$ws.Range("I34").Value = 2035.7142
$ws.Range("J34").Value = 4165.68
$ws.Range("K34").Value = 2035.7142
$ws.Range("L34").Value = 4165.68
$ws.Range("M34").Value = -1833.7142
$ws.Range("N34").Value = -4569.68
$ws.Range("H45").Value = 6037
$ws.Range("I45").Value = 3000
$ws.Range("J45").Value = 9074
$ws.Range("K45").Value = 3000
$ws.Range("L45").Value = 9074
$ws.Range("M45").Value = -2407
$ws.Range("N45").Value = -10260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 326.6154
$ws.Range("I14").Value = 326.6154
$ws.Range("K14").Value = 979.8462000000001
$ws.Range("M14").Value = -806.8462000000001
$ws.Range("H81").Value = 123751624
$ws.Range("I81").Value = 499.5
$ws.Range("J81").Value = 165002000
$ws.Range("K81").Value = 1498.5
$ws.Range("L81").Value = 495006000
$ws.Range("M81").Value = -375.5
$ws.Range("N81").Value = -495008246
$ws.Range("H84").Value = 123751624
$ws.Range("I84").Value = 499.5
$ws.Range("J84").Value = 165002000
$ws.Range("K84").Value = 4495.5
$ws.Range("L84").Value = 1485018000
$ws.Range("M84").Value = 1120.5
$ws.Range("N84").Value = -1485029232
$ws.Range("H107").Value = 300820.38
$ws.Range("I107").Value = 471.3913
$ws.Range("J107").Value = 928822.8
$ws.Range("K107").Value = 1414.1739
$ws.Range("L107").Value = 2786468.4
$ws.Range("M107").Value = 505.8261
$ws.Range("N107").Value = -2790308.4
$ws.Range("H113").Value = 351763.8
$ws.Range("I113").Value = 605
$ws.Range("K113").Value = 1815
$ws.Range("M113").Value = 355
$ws.Range("H131").Value = 877.38
$ws.Range("J131").Value = 945.55817
$ws.Range("L131").Value = 2836.67451
$ws.Range("N131").Value = -12916.67451
$ws.Range("H132").Value = 602831.9399999999
$ws.Range("I132").Value = 1197518
$ws.Range("J132").Value = 8145.909
$ws.Range("K132").Value = 10777662
$ws.Range("L132").Value = 73313.181
$ws.Range("M132").Value = -10775132
$ws.Range("N132").Value = -78373.181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 26329
$ws.Range("J62").Value = 26329
$ws.Range("L62").Value = 26329
$ws.Range("N62").Value = -27701
$ws.Range("H65").Value = 26329
$ws.Range("J65").Value = 26329
$ws.Range("L65").Value = 78987
$ws.Range("N65").Value = -85851
$ws.Range("H70").Value = 47221.543
$ws.Range("I70").Value = 55700.4
$ws.Range("J70").Value = 4827.25
$ws.Range("K70").Value = 55700.4
$ws.Range("L70").Value = 4827.25
$ws.Range("M70").Value = -55430.4
$ws.Range("N70").Value = -5367.25
$ws.Range("H73").Value = 47221.543
$ws.Range("I73").Value = 55700.4
$ws.Range("J73").Value = 4827.25
$ws.Range("K73").Value = 55700.4
$ws.Range("L73").Value = 4827.25
$ws.Range("M73").Value = -54764.4
$ws.Range("N73").Value = -6699.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 18046
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 18046
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 18046
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -18860
$ws.Range("H48").Value = 16029
$ws.Range("I48").Value = 10020.5
$ws.Range("K48").Value = 10020.5
$ws.Range("M48").Value = -9359.5
$ws.Range("H63").Value = 32988.332
$ws.Range("J63").Value = 32988.332
$ws.Range("L63").Value = 32988.332
$ws.Range("N63").Value = -34486.332
$ws.Range("H66").Value = 32988.332
$ws.Range("J66").Value = 32988.332
$ws.Range("L66").Value = 98964.99600000001
$ws.Range("N66").Value = -106452.996
$ws.Range("H68").Value = 3074.875
$ws.Range("I68").Value = 1400
$ws.Range("J68").Value = 4749.75
$ws.Range("K68").Value = 1400
$ws.Range("L68").Value = 4749.75
$ws.Range("M68").Value = -651
$ws.Range("N68").Value = -6247.75
$ws.Range("H71").Value = 3074.875
$ws.Range("I71").Value = 1400
$ws.Range("J71").Value = 4749.75
$ws.Range("K71").Value = 7000
$ws.Range("L71").Value = 23748.75
$ws.Range("M71").Value = -3256
$ws.Range("N71").Value = -31236.75
$ws.Range("H132").Value = 4172.892
$ws.Range("I132").Value = 4330.8696
$ws.Range("J132").Value = 3913.3572
$ws.Range("K132").Value = 12992.6088
$ws.Range("L132").Value = 11740.0716
$ws.Range("M132").Value = -10462.6088
$ws.Range("N132").Value = -16800.0716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9000
$ws.Range("J41").Value = 9000
$ws.Range("L41").Value = 9000
$ws.Range("N41").Value = -9780
$ws.Range("H62").Value = 9156.25
$ws.Range("I62").Value = 3837.5
$ws.Range("J62").Value = 14475
$ws.Range("K62").Value = 3837.5
$ws.Range("L62").Value = 14475
$ws.Range("M62").Value = -3213.5
$ws.Range("N62").Value = -15723
$ws.Range("H65").Value = 9156.25
$ws.Range("I65").Value = 3837.5
$ws.Range("J65").Value = 14475
$ws.Range("K65").Value = 19187.5
$ws.Range("L65").Value = 72375
$ws.Range("M65").Value = -16067.5
$ws.Range("N65").Value = -78615
$ws.Range("H125").Value = 27263.637
$ws.Range("J125").Value = 27263.637
$ws.Range("L125").Value = 27263.637
$ws.Range("N125").Value = -37103.637

